$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new normalization rule rows after the existing data (rows 1-13 already populated)
$ws.Range("A14").Value = "p(i|e+)"
$ws.Range("B14").Value = "p"

$ws.Range("A15").Value = "di"
$ws.Range("B15").Value = "d"

$ws.Range("A16").Value = "bi"
$ws.Range("B16").Value = "b"
$ws.Range("C16").Value = "bichwala"

# Match the row heights used by the rest of the sheet and leave a few
# blank rows below the new data (same as the rest of the worksheet).
$ws.Rows(14).RowHeight = 12.1
$ws.Rows(15).RowHeight = 12.1
$ws.Rows(16).RowHeight = 12.8
$ws.Rows(17).RowHeight = 12.1
$ws.Rows(18).RowHeight = 12.1
$ws.Rows(19).RowHeight = 12.1

# Move the selection to the next empty row, as in the saved workbook
$ws.Range("A17:C17").Select()
